# Weekly update: insert a new price record (row 96) for Camote at
# "Vega Modelo de Temuco" and push the existing rows 96-120 down to 97-121.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 96 (shifts 96..120 -> 97..121)
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new weekly record
$ws.Range("A96").Value2 = 10
$ws.Range("B96").Value2 = "Vega Modelo de Temuco"
$ws.Range("C96").Value2 = "La Araucanía"
$ws.Range("D96").Value2 = 44809
$ws.Range("E96").Value2 = 9
$ws.Range("F96").Value2 = 100114002
$ws.Range("G96").Value2 = "Camote"
$ws.Range("H96").Value2 = "Sin especificar"
$ws.Range("I96").Value2 = "Primera"
$ws.Range("J96").Value2 = 50
$ws.Range("K96").Value2 = 20000
$ws.Range("L96").Value2 = 20000
$ws.Range("M96").Value2 = 20000
$ws.Range("N96").Value2 = "`$/malla 20 kilos"
$ws.Range("O96").Value2 = "Perú"
$ws.Range("P96").Value2 = 1000
$ws.Range("Q96").Value2 = 20
$ws.Range("R96").Value2 = "Hortaliza"
